$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1949367088607595
$ws.Range("C2").Value = 0.5468354430379747
$ws.Range("J2").Value = 0.02025316455696203
$ws.Range("P2").Value = 0.1518987341772152
$ws.Range("S2").Value = 0.08607594936708861
$ws.Range("B3").Value = 0.01777777777777778
$ws.Range("C3").Value = 0.04
$ws.Range("J3").Value = 0.04
$ws.Range("P3").Value = 0.6888888888888889
$ws.Range("S3").Value = 0.2133333333333333
$ws.Range("J4").Value = 0.1224489795918367
$ws.Range("P4").Value = 0.673469387755102
$ws.Range("S4").Value = 0.2040816326530612
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.0846774193548387
$ws.Range("D6").Value = 0.01209677419354839
$ws.Range("F6").Value = 0.0282258064516129
$ws.Range("J6").Value = 0.2862903225806452
$ws.Range("O6").Value = 0.004032258064516129
$ws.Range("Q6").Value = 0.1653225806451613
$ws.Range("R6").Value = 0.07661290322580645
$ws.Range("S6").Value = 0.3427419354838709
$ws.Range("B7").Value = 0.1684210526315789
$ws.Range("D7").Value = 0.02105263157894737
$ws.Range("F7").Value = 0.03684210526315789
$ws.Range("J7").Value = 0.1947368421052632
$ws.Range("O7").Value = 0.02105263157894737
$ws.Range("Q7").Value = 0.1368421052631579
$ws.Range("R7").Value = 0.09473684210526316
$ws.Range("S7").Value = 0.3263157894736842
$ws.Range("B8").Value = 0.1111111111111111
$ws.Range("D8").Value = 0.01694915254237288
$ws.Range("F8").Value = 0.04896421845574388
$ws.Range("J8").Value = 0.1261770244821092
$ws.Range("O8").Value = 0.01506591337099812
$ws.Range("Q8").Value = 0.199623352165725
$ws.Range("R8").Value = 0.07721280602636535
$ws.Range("S8").Value = 0.4048964218455744
$ws.Range("B9").Value = 0.08144796380090498
$ws.Range("D9").Value = 0.009049773755656109
$ws.Range("E9").Value = 0.004524886877828055
$ws.Range("F9").Value = 0.08144796380090498
$ws.Range("J9").Value = 0.1493212669683258
$ws.Range("O9").Value = 0.01809954751131222
$ws.Range("Q9").Value = 0.1719457013574661
$ws.Range("R9").Value = 0.08144796380090498
$ws.Range("S9").Value = 0.4027149321266968
$ws.Range("B10").Value = 0.1136662286465177
$ws.Range("D10").Value = 0.019053876478318
$ws.Range("F10").Value = 0.07161629434954007
$ws.Range("J10").Value = 0.1103810775295664
$ws.Range("O10").Value = 0.01839684625492773
$ws.Range("Q10").Value = 0.2240473061760841
$ws.Range("R10").Value = 0.06504599211563732
$ws.Range("S10").Value = 0.3777923784494087
$ws.Range("G11").Value = 0.1889250814332248
$ws.Range("J11").Value = 0.1009771986970684
$ws.Range("K11").Value = 0.244299674267101
$ws.Range("L11").Value = 0.4527687296416938
$ws.Range("S11").Value = 0.01302931596091205
$ws.Range("G12").Value = 0.8043478260869565
$ws.Range("J12").Value = 0.1521739130434783
$ws.Range("L12").Value = 0.007246376811594203
$ws.Range("S12").Value = 0.03623188405797102
$ws.Range("G13").Value = 0.5777777777777777
$ws.Range("J13").Value = 0.3777777777777778
$ws.Range("S13").Value = 0.04444444444444445
$ws.Range("F15").Value = 0.007843137254901961
$ws.Range("H15").Value = 0.207843137254902
$ws.Range("I15").Value = 0.07058823529411765
$ws.Range("J15").Value = 0.3686274509803922
$ws.Range("K15").Value = 0.05098039215686274
$ws.Range("M15").Value = 0.007843137254901961
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.2196078431372549
$ws.Range("F16").Value = 0.01694915254237288
$ws.Range("H16").Value = 0.1779661016949153
$ws.Range("I16").Value = 0.09322033898305085
$ws.Range("J16").Value = 0.4067796610169492
$ws.Range("K16").Value = 0.1016949152542373
$ws.Range("M16").Value = 0.01271186440677966
$ws.Range("O16").Value = 0.06779661016949153
$ws.Range("S16").Value = 0.1228813559322034
$ws.Range("F17").Value = 0.01824817518248175
$ws.Range("H17").Value = 0.1733576642335767
$ws.Range("I17").Value = 0.08029197080291971
$ws.Range("J17").Value = 0.4197080291970803
$ws.Range("K17").Value = 0.0948905109489051
$ws.Range("M17").Value = 0.01824817518248175
$ws.Range("N17").Value = 0.001824817518248175
$ws.Range("O17").Value = 0.09124087591240876
$ws.Range("S17").Value = 0.1021897810218978
$ws.Range("F18").Value = 0.02061855670103093
$ws.Range("H18").Value = 0.1855670103092784
$ws.Range("I18").Value = 0.07731958762886598
$ws.Range("J18").Value = 0.4639175257731959
$ws.Range("K18").Value = 0.1082474226804124
$ws.Range("O18").Value = 0.05154639175257732
$ws.Range("S18").Value = 0.09278350515463918
$ws.Range("F19").Value = 0.01639344262295082
$ws.Range("H19").Value = 0.2230933713471133
$ws.Range("I19").Value = 0.08624376336421953
$ws.Range("J19").Value = 0.3955808980755524
$ws.Range("K19").Value = 0.08766928011404133
$ws.Range("M19").Value = 0.02209550962223806
$ws.Range("N19").Value = 0.0007127583749109052
$ws.Range("O19").Value = 0.0684248039914469
$ws.Range("S19").Value = 0.09978617248752673
